$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, pushing existing rows 11.. down by one.
$ws.Rows("11").Insert()

# Fill in the newly inserted row 11 with the new data point.
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Vega Modelo de Temuco"
$ws.Range("C11").Value = "La Araucanía"
$ws.Range("D11").Value = 45111
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = 100112010
$ws.Range("G11").Value = "Achicoria"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 8500
$ws.Range("L11").Value = 8500
$ws.Range("M11").Value = 8500
$ws.Range("N11").Value = "`$/caja 18 unidades"
$ws.Range("O11").Value = "Región Metropolitana"
$ws.Range("P11").Value = 472
$ws.Range("Q11").Value = 18
$ws.Range("R11").Value = "Hortaliza"
